$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (cell reference, new value). All of these cells already hold
# text in the workbook (coin prices/volumes are stored as strings, not
# numbers). Several of the new price values are plain decimals (e.g.
# "230.91"); a leading apostrophe forces Excel to keep them as text (quote
# -prefixed), matching the original t="inlineStr" cells instead of letting
# COM auto-convert them to numeric values.
$changes = @(
    ,@('D2', '44.121.03')
    ,@('E2', '  +5.34%  ')
    ,@('D3', '2.294.06')
    ,@('E3', '  +2.69%  ')
    ,@('E4', '  -0.02%  ')
    ,@('D5', '''230.91')
    ,@('E5', '  -0.75%  ')
    ,@('D6', '''0.628')
    ,@('E6', '  -0.03%  ')
    ,@('D7', '''61.04')
    ,@('E7', '  -0.60%  ')
    ,@('E8', '  -0.09%  ')
    ,@('D9', '''0.425')
    ,@('E9', '  +4.95%  ')
    ,@('D10', '''0.0947')
    ,@('E10', '  +4.10%  ')
    ,@('E11', '  +0.22%  ')
    ,@('D12', '2.627.08')
    ,@('E12', '  +2.42%  ')
    ,@('D13', '''24.40')
    ,@('E13', '  +8.94%  ')
    ,@('D14', '''15.80')
    ,@('E14', '  +0.68%  ')
    ,@('D15', '''5.88')
    ,@('E15', '  +4.85%  ')
    ,@('D16', '''0.815')
    ,@('E16', '  +1.29%  ')
    ,@('D17', '2.299.18')
    ,@('E17', '  +2.65%  ')
    ,@('D18', '44.080.27')
    ,@('E18', '  +5.27%  ')
    ,@('D19', '0.0₃0947')
    ,@('E19', '  +4.16%  ')
    ,@('D20', '''73.59')
    ,@('E20', '  +1.50%  ')
    ,@('D21', '''6.27')
    ,@('E21', '  +3.86%  ')
    ,@('D22', '''254.81')
    ,@('E22', '  +1.08%  ')
    ,@('E23', '  -0.05%  ')
    ,@('D24', '''2.55')
    ,@('E24', '  +6.42%  ')
    ,@('E25', '  -0.52%  ')
    ,@('D26', '''9.93')
    ,@('E26', '  +2.04%  ')
    ,@('D27', '''171.26')
    ,@('E27', '  +1.24%  ')
    ,@('E28', '  -2.32%  ')
    ,@('D29', '''20.63')
    ,@('E29', '  +2.73%  ')
    ,@('D30', '''1.43')
    ,@('E30', '  +0.10%  ')
    ,@('D31', '''2.76')
    ,@('E31', '  +1.31%  ')
    ,@('E32', '  +0.02%  ')
    ,@('D33', '''5.08')
    ,@('E33', '  +0.19%  ')
    ,@('D34', '''4.79')
    ,@('E34', '  +2.47%  ')
    ,@('D35', '''0.0660')
    ,@('E35', '  +3.34%  ')
    ,@('D36', '''6.52')
    ,@('E36', '  -2.04%  ')
    ,@('D37', '''2.41')
    ,@('E37', '  +1.54%  ')
    ,@('D38', '''3.64')
    ,@('E38', '  -2.81%  ')
    ,@('E39', '  +4.24%  ')
    ,@('E40', '  -0.05%  ')
    ,@('D41', '''8.87')
    ,@('E41', '  +3.87%  ')
    ,@('D42', '''0.000225')
    ,@('E42', '  -12.93%  ')
    ,@('D43', '''0.0973')
    ,@('E43', '  +1.15%  ')
    ,@('D44', '''99.09')
    ,@('E44', '  -0.48%  ')
    ,@('D45', '''1.22')
    ,@('E45', '  -1.44%  ')
    ,@('B46', 'Maker')
    ,@('C46', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr')
    ,@('D46', '1.481.96')
    ,@('E46', '  -0.01%  ')
    ,@('B47', 'InjectiveProtocol')
    ,@('C47', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj')
    ,@('D47', '''16.97')
    ,@('E47', '  +2.03%  ')
    ,@('B48', 'FTXToken')
    ,@('C48', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt')
    ,@('D48', '''4.39')
    ,@('E48', '  -7.77%  ')
    ,@('B49', 'Celestia')
    ,@('C49', 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia')
    ,@('D49', '''10.17')
    ,@('E49', '  +15.13%  ')
    ,@('D50', '''1.10')
    ,@('E50', '  +1.62%  ')
    ,@('D51', '''2.27')
    ,@('E51', '  +5.94%  ')
)

foreach ($change in $changes) {
    $ref = $change[0]
    $val = $change[1]
    $ws.Range($ref).Value = $val
}

Write-Host "Applied $($changes.Count) cell updates"
